$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 12:50"

# España (row 7) - refreshed case counts
$ws.Range("B7").Value = 78797
$ws.Range("C7").Value = 5562
$ws.Range("D7").Value = 14709
$ws.Range("E7").Value = 57560
$ws.Range("F7").Value = 4165
$ws.Range("G7").Value = 546
$ws.Range("H7").Value = 6528

# Barein (row 63) - refreshed case counts
$ws.Range("B63").Value = 499
$ws.Range("C63").Value = 23
$ws.Range("D63").Value = 272
$ws.Range("E63").Value = 223

# Senegal's case count overtakes Ghana/Costa de Marfil/Uzbekistan, so the
# sorted country list shuffles rows 96-99 (row 100, Brunei, is untouched).
$ws.Range("A96").Value = "Senegal"
$ws.Range("B96").Value = 142
$ws.Range("C96").Value = 12
$ws.Range("D96").Value = 27
$ws.Range("E96").Value = 115
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0

$ws.Range("A97").Value = "Ghana"
$ws.Range("B97").Value = 141
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 2
$ws.Range("E97").Value = 134
$ws.Range("F97").Value = 1
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 5

$ws.Range("A98").Value = "Costa de Marfil"
$ws.Range("B98").Value = 140
$ws.Range("C98").Value = 39
$ws.Range("D98").Value = 3
$ws.Range("E98").Value = 137
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0

$ws.Range("A99").Value = "Uzbekistan"
$ws.Range("B99").Value = 133
$ws.Range("C99").Value = 29
$ws.Range("D99").Value = 7
$ws.Range("E99").Value = 124
$ws.Range("F99").Value = 8
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 2
